$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Calibration data rows (A2:D12) re-sorted chronologically (column A, time,
# ascending) after performing calibration of the needle.
$data = @(
    @(53613.535026, -0.000023758669731, -0.000019758389067, -0.0000079040205007),
    @(53624.335026, -0.00017156735296, -0.00014172239333, -0.000053008852357),
    @(53635.935027, -0.0003666249, -0.0003011899, -0.0001014434),
    @(53645.803027, -0.0005438759, -0.0004508612, -0.0001518111),
    @(53656.199028, -0.0007193427, -0.0006048907, -0.0002048983),
    @(53666.935028, -0.0008884488, -0.0007603567, -0.0002556789),
    @(53690.999029, -0.0007062664, -0.0005961034, -0.0002109604),
    @(53702.19903,  -0.0005284251, -0.000440739,  -0.0001587659),
    @(53712.33503,  -0.0003495365, -0.0002892978, -0.0001055855),
    @(53723.667031, -0.00015841520936, -0.00013209913367, -0.000053904856423),
    @(53734.531031, -0.000025406350801, -0.00002165626666, -0.000010544177038)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
